$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the same style as the existing header cells (e.g. A1) to the new
# header range before writing values, so formatting (bold, border, centered)
# matches the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# New header cells F1:H1
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Boolean data for rows 2-23 in columns F (KNN), G (SVM), H (RF)
# Default to FALSE, with TRUE flags for the outlier rows (9 and 18)
$outlierRows = @(9, 18)

for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 6).Value = $false
    if ($outlierRows -contains $r) {
        $ws.Cells.Item($r, 7).Value = $true
        $ws.Cells.Item($r, 8).Value = $true
    } else {
        $ws.Cells.Item($r, 7).Value = $false
        $ws.Cells.Item($r, 8).Value = $false
    }
}
